# Realistic Polymarket fills with order book walking, skip trades with high slippage
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: update aggregate performance metrics
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1000.11   # Current Capital
$wsSummary.Range("B4").Value = 0.12      # Total P&L $
$wsSummary.Range("B5").Value = 0.22      # Total P&L %
$wsSummary.Range("B6").Value = 11        # Total Trades
$wsSummary.Range("B8").Value = 4         # Losing Trades
$wsSummary.Range("B9").Value = 45.45     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: update the "leadlag" strategy row (row 5)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 100.09   # Capital
$wsStatus.Range("D5").Value = 10       # Trades
$wsStatus.Range("E5").Value = 0.1      # P&L $
$wsStatus.Range("F5").Value = 0.09     # P&L %
$wsStatus.Range("G5").Value = 40       # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet: trade #11 (row 12) moves from OPEN to CLOSED with a
# realistic fill that walked the order book
# ---------------------------------------------------------------------------
$wsTrades = $wb.Worksheets.Item("All Trades")
$wsTrades.Range("G12").Value = 67860.523502   # Exit Price
$wsTrades.Range("H12").Value = "CLOSED"       # Status
$wsTrades.Range("I12").Value = -0.7048        # P&L %
$wsTrades.Range("J12").Value = -0.04          # P&L $
$wsTrades.Range("K12").Value = 100.09         # Capital After
$wsTrades.Range("N12").Value = "time_exit_5min"  # Exit Reason
$wsTrades.Range("O12").Value = 5              # Minutes Held

# ---------------------------------------------------------------------------
# leadlag sheet: same trade, mirrored on the strategy-specific sheet (row 11)
# ---------------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")
$wsLeadlag.Range("G11").Value = 67860.523502   # Exit Price
$wsLeadlag.Range("H11").Value = "CLOSED"       # Status
$wsLeadlag.Range("I11").Value = -0.7048        # P&L %
$wsLeadlag.Range("J11").Value = -0.04          # P&L $
$wsLeadlag.Range("K11").Value = 100.09         # Capital After
$wsLeadlag.Range("N11").Value = "time_exit_5min"  # Exit Reason
$wsLeadlag.Range("O11").Value = 5              # Minutes Held
